$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Content fixes: the "Observaciones"/"Pagina en ficha tecnica" text
# for the "Disponibilidad de lideres en piso" requirement (row 10)
# had been entered one row too low (row 11). Move it up to row 10
# and replace row 11's values with the correct text for that row.
# ------------------------------------------------------------------

# Row 10 ("Disponibilidad de lideres en piso..."): fill in D10/E10
$ws.Range("D10").Value = "Ver plan en la evaluacion técnica (Pág. X)."
$ws.Range("E10").Value = "La disponibilidad de líderes garantiza una reducción del 95% en errores de armado y etiquetado, y una respuesta inmediata ante variaciones en la carga de contenedores."

# Row 11 ("Entrega de reportes diarios/semanales"): clear D11, set new E11
$ws.Range("D11").ClearContents()
$ws.Range("E11").Value = "Transparencia total de datos."

# ------------------------------------------------------------------
# Re-apply the formatting (wrap text, vertical-center alignment) on
# the touched cells so they keep the same look as the rest of the
# "Observaciones"/"Pagina en ficha tecnica" columns.
# ------------------------------------------------------------------
$touched = @("D10", "E10", "D11")
foreach ($addr in $touched) {
    $c = $ws.Range($addr)
    $c.WrapText = $true
    $c.VerticalAlignment = -4108
    $c.ShrinkToFit = $false
}

# ------------------------------------------------------------------
# Normalize formatting on the "Observaciones" column (E2:E14, except
# E11 which already keeps its formatting) so it matches the rest of
# the sheet.
# ------------------------------------------------------------------
$obsRows = @(2,3,4,5,6,7,8,9,12,13,14)
foreach ($r in $obsRows) {
    $c = $ws.Cells.Item($r, 5)
    $c.WrapText = $true
    $c.VerticalAlignment = -4108
    $c.ShrinkToFit = $false
}

# ------------------------------------------------------------------
# Auto-fit rows 5-14 so the explicit custom row height is cleared and
# rows go back to their natural auto-computed height.
# ------------------------------------------------------------------
$ws.Rows("5:14").AutoFit()
